$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 is the student "Полторабатько Кирилл" (student index 24) - fill in
# the homework scores (columns D..L) that were previously left blank.

# D27:I27 already carry the "filled" style (s=2); just populate the values.
$ws.Range("D27:I27").Value = 5

# J27 and K27 need both a value and the "filled" styling used elsewhere in
# the sheet for those columns (borders/fill differ per column group), so
# copy formatting from already-styled neighbor cells before setting values.
$ws.Range("J24").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("J27").Value = 5

$ws.Range("K12").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K27").Value = 5

$ws.Range("L12").Copy()
$ws.Range("L27").PasteSpecial(-4122)
$ws.Range("L27").Value = 10

$excel.CutCopyMode = 0

# Update the active selection to match the new state (T27 instead of T28)
$ws.Range("T27").Select()
